$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 39-52 (re-sequenced dates/quality/volume/prices) ---
$ws.Range("D39").Value = 44455
$ws.Range("L39").Value = "Especial"
$ws.Range("M39").Value = 15
$ws.Range("N39").Value = 3500
$ws.Range("O39").Value = 3500
$ws.Range("P39").Value = 3500
$ws.Range("Q39").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S39").Value = 3500
$ws.Range("T39").Value = 1

$ws.Range("D40").Value = 44455
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 20
$ws.Range("N40").Value = 3000
$ws.Range("O40").Value = 3000
$ws.Range("P40").Value = 3000
$ws.Range("Q40").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S40").Value = 3000
$ws.Range("T40").Value = 1

$ws.Range("D41").Value = 44414
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 55
$ws.Range("N41").Value = 3500
$ws.Range("O41").Value = 3500
$ws.Range("P41").Value = 3500
$ws.Range("Q41").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S41").Value = 3500
$ws.Range("T41").Value = 1

$ws.Range("D42").Value = 44165
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 50
$ws.Range("N42").Value = 2300
$ws.Range("O42").Value = 2300
$ws.Range("P42").Value = 2300
$ws.Range("Q42").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S42").Value = 2300
$ws.Range("T42").Value = 1

$ws.Range("D43").Value = 44427
$ws.Range("L43").Value = "Especial"
$ws.Range("M43").Value = 65
$ws.Range("N43").Value = 24000
$ws.Range("O43").Value = 24000
$ws.Range("P43").Value = 24000
$ws.Range("Q43").Value = "$/bandeja 7 kilos"
$ws.Range("S43").Value = 3429
$ws.Range("T43").Value = 7

$ws.Range("D44").Value = 44447
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 40
$ws.Range("N44").Value = 3000
$ws.Range("O44").Value = 3000
$ws.Range("P44").Value = 3000
$ws.Range("Q44").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S44").Value = 3000
$ws.Range("T44").Value = 1

$ws.Range("D45").Value = 44187
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 15
$ws.Range("N45").Value = 3200
$ws.Range("O45").Value = 3200
$ws.Range("P45").Value = 3200
$ws.Range("Q45").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S45").Value = 3200
$ws.Range("T45").Value = 1

$ws.Range("D46").Value = 44433
$ws.Range("L46").Value = "Especial"
$ws.Range("M46").Value = 20
$ws.Range("N46").Value = 4500
$ws.Range("O46").Value = 4500
$ws.Range("P46").Value = 4500
$ws.Range("Q46").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S46").Value = 4500
$ws.Range("T46").Value = 1

$ws.Range("D47").Value = 44438
$ws.Range("L47").Value = "Especial"
$ws.Range("M47").Value = 35
$ws.Range("N47").Value = 3500
$ws.Range("O47").Value = 3500
$ws.Range("P47").Value = 3500
$ws.Range("Q47").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S47").Value = 3500
$ws.Range("T47").Value = 1

$ws.Range("D48").Value = 44438
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 20
$ws.Range("N48").Value = 3000
$ws.Range("O48").Value = 3000
$ws.Range("P48").Value = 3000
$ws.Range("Q48").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S48").Value = 3000
$ws.Range("T48").Value = 1

$ws.Range("D49").Value = 44160
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 120
$ws.Range("N49").Value = 2200
$ws.Range("O49").Value = 2300
$ws.Range("P49").Value = 2246
$ws.Range("Q49").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S49").Value = 2246
$ws.Range("T49").Value = 1

$ws.Range("D50").Value = 44162
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 85
$ws.Range("N50").Value = 2200
$ws.Range("O50").Value = 2300
$ws.Range("P50").Value = 2247
$ws.Range("Q50").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S50").Value = 2247
$ws.Range("T50").Value = 1

$ws.Range("D51").Value = 44411
$ws.Range("L51").Value = "Segunda"
$ws.Range("M51").Value = 10
$ws.Range("N51").Value = 3000
$ws.Range("O51").Value = 3000
$ws.Range("P51").Value = 3000
$ws.Range("Q51").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S51").Value = 3000
$ws.Range("T51").Value = 1

$ws.Range("D52").Value = 44425
$ws.Range("L52").Value = "Especial"
$ws.Range("M52").Value = 35
$ws.Range("N52").Value = 4500
$ws.Range("O52").Value = 4500
$ws.Range("P52").Value = 4500
$ws.Range("Q52").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S52").Value = 4500
$ws.Range("T52").Value = 1

# --- Add new rows 53-54 ---
$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "Vega Modelo de Temuco"
$ws.Range("C53").Value = "La Araucanía"
$ws.Range("D53").Value = 44425
$ws.Range("D53").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E53").Value = 9
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100107
$ws.Range("H53").Value = "Otros"
$ws.Range("I53").Value = 100107002
$ws.Range("J53").Value = "Chirimoya"
$ws.Range("K53").Value = "Cultivar IV Región"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 20
$ws.Range("N53").Value = 3500
$ws.Range("O53").Value = 3500
$ws.Range("P53").Value = 3500
$ws.Range("Q53").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R53").Value = "Provincia del Elquí"
$ws.Range("S53").Value = 3500
$ws.Range("T53").Value = 1

$ws.Range("A54").Value = 10
$ws.Range("B54").Value = "Vega Modelo de Temuco"
$ws.Range("C54").Value = "La Araucanía"
$ws.Range("D54").Value = 44425
$ws.Range("D54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100107
$ws.Range("H54").Value = "Otros"
$ws.Range("I54").Value = 100107002
$ws.Range("J54").Value = "Chirimoya"
$ws.Range("K54").Value = "Cultivar IV Región"
$ws.Range("L54").Value = "Segunda"
$ws.Range("M54").Value = 25
$ws.Range("N54").Value = 3000
$ws.Range("O54").Value = 3000
$ws.Range("P54").Value = 3000
$ws.Range("Q54").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R54").Value = "Provincia del Elquí"
$ws.Range("S54").Value = 3000
$ws.Range("T54").Value = 1
